$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2750
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -270
$ws.Range("N20").Value = -5460

$ws.Range("H33").Value = 141.83333
$ws.Range("I33").Value = 127.454544
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 127.454544
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = 101.545456
$ws.Range("N33").Value = -758

$ws.Range("H35").Value = 2750
$ws.Range("I35").Value = 500
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 500
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -121
$ws.Range("N35").Value = -5758

$ws.Range("H43").Value = 11457.8
$ws.Range("I43").Value = 12099.667
$ws.Range("K43").Value = 12099.667
$ws.Range("M43").Value = -12030.667

$ws.Range("H62").Value = 4224.8335
$ws.Range("I62").Value = 4737.5
$ws.Range("K62").Value = 4737.5
$ws.Range("M62").Value = -4113.5

$ws.Range("H65").Value = 4224.8335
$ws.Range("I65").Value = 4737.5
$ws.Range("K65").Value = 23687.5
$ws.Range("M65").Value = -20567.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2726.6897
$ws.Range("I32").Value = 1623.1538
$ws.Range("K32").Value = 1623.1538
$ws.Range("M32").Value = -1336.1538

$ws.Range("H34").Value = 37499.5

$ws.Range("H45").Value = 2401.9167
$ws.Range("I45").Value = 2302.5557
$ws.Range("J45").Value = 2700
$ws.Range("K45").Value = 2302.5557
$ws.Range("L45").Value = 2700
$ws.Range("M45").Value = -1925.5557
$ws.Range("N45").Value = -3454

$ws.Range("H54").Value = 88999.5
$ws.Range("J54").Value = 88999.5
$ws.Range("L54").Value = 88999.5
$ws.Range("N54").Value = -90537.5

$ws.Range("H61").Value = 1973
$ws.Range("I61").Value = 1598.3846
$ws.Range("K61").Value = 1598.3846
$ws.Range("M61").Value = -1386.3846

$ws.Range("H88").Value = 3250.5386
$ws.Range("I88").Value = 2799.75
$ws.Range("J88").Value = 3450.889
$ws.Range("K88").Value = 2799.75
$ws.Range("L88").Value = 3450.889
$ws.Range("M88").Value = -2393.75
$ws.Range("N88").Value = -4262.889

$ws.Range("H91").Value = 3250.5386
$ws.Range("I91").Value = 2799.75
$ws.Range("J91").Value = 3450.889
$ws.Range("K91").Value = 2799.75
$ws.Range("L91").Value = 3450.889
$ws.Range("M91").Value = -1395.75
$ws.Range("N91").Value = -6258.889

$ws.Range("H97").Value = 432.13043
$ws.Range("I97").Value = 274.6111
$ws.Range("J97").Value = 999.2
$ws.Range("K97").Value = 274.6111
$ws.Range("L97").Value = 999.2
$ws.Range("M97").Value = 221.3889
$ws.Range("N97").Value = -1991.2

$ws.Range("H136").Value = 1973
$ws.Range("I136").Value = 1598.3846
$ws.Range("K136").Value = 4795.1538
$ws.Range("M136").Value = -2245.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1338.75
$ws.Range("I64").Value = 1200
$ws.Range("J64").Value = 1385
$ws.Range("K64").Value = 1200
$ws.Range("L64").Value = 1385
$ws.Range("M64").Value = -975
$ws.Range("N64").Value = -1835

$ws.Range("H67").Value = 1338.75
$ws.Range("I67").Value = 1200
$ws.Range("J67").Value = 1385
$ws.Range("K67").Value = 1200
$ws.Range("L67").Value = 1385
$ws.Range("M67").Value = -420
$ws.Range("N67").Value = -2945

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 3293.1365
$ws.Range("J94").Value = 4622.5557
$ws.Range("L94").Value = 4622.5557
$ws.Range("N94").Value = -5524.5557

$ws.Range("H134").Value = 2601.18
$ws.Range("I134").Value = 2436.9512
$ws.Range("K134").Value = 7310.8536
$ws.Range("M134").Value = -4775.8536

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 39999.668
$ws.Range("I13").Value = 15000
$ws.Range("J13").Value = 52499.5
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 52499.5
$ws.Range("M13").Value = -14861
$ws.Range("N13").Value = -52777.5

$ws.Range("H31").Value = 1845.5294
$ws.Range("I31").Value = 1713.375
$ws.Range("J31").Value = 1963
$ws.Range("K31").Value = 1713.375
$ws.Range("L31").Value = 1963
$ws.Range("M31").Value = -1418.375
$ws.Range("N31").Value = -2553

$ws.Range("H34").Value = 1845.5294
$ws.Range("I34").Value = 1713.375
$ws.Range("J34").Value = 1963
$ws.Range("K34").Value = 1713.375
$ws.Range("L34").Value = 1963
$ws.Range("M34").Value = -1511.375
$ws.Range("N34").Value = -2367

$ws.Range("H58").Value = 2560.2273
$ws.Range("I58").Value = 2671.9167
$ws.Range("J58").Value = 2426.2
$ws.Range("K58").Value = 2671.9167
$ws.Range("L58").Value = 2426.2
$ws.Range("M58").Value = -2468.9167
$ws.Range("N58").Value = -2832.2

$ws.Range("H107").Value = 993.94116
$ws.Range("I107").Value = 747.75
$ws.Range("J107").Value = 1026.7667
$ws.Range("K107").Value = 747.75
$ws.Range("L107").Value = 1026.7667
$ws.Range("M107").Value = 1172.25
$ws.Range("N107").Value = -4866.7667

$ws.Range("H136").Value = 2560.2273
$ws.Range("I136").Value = 2671.9167
$ws.Range("J136").Value = 2426.2
$ws.Range("K136").Value = 8015.750100000001
$ws.Range("L136").Value = 7278.599999999999
$ws.Range("M136").Value = -5465.750100000001
$ws.Range("N136").Value = -12378.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 83343990
$ws.Range("I4").Value = 125000616
$ws.Range("K4").Value = 375001848
$ws.Range("M4").Value = -375001736

$ws.Range("H8").Value = 608.46155
$ws.Range("I8").Value = 608.46155
$ws.Range("K8").Value = 1825.38465
$ws.Range("M8").Value = -1686.38465

$ws.Range("H107").Value = 772.9167
$ws.Range("I107").Value = 688
$ws.Range("J107").Value = 833.5714
$ws.Range("K107").Value = 2064
$ws.Range("L107").Value = 2500.7142
$ws.Range("M107").Value = -144
$ws.Range("N107").Value = -6340.7142

$ws.Range("H109").Value = 3799
$ws.Range("I109").Value = 3133.111
$ws.Range("J109").Value = 4997.6
$ws.Range("K109").Value = 9399.332999999999
$ws.Range("L109").Value = 14992.8
$ws.Range("M109").Value = -8359.332999999999
$ws.Range("N109").Value = -17072.8

$ws.Range("H131").Value = 1377.8125
$ws.Range("I131").Value = 869
$ws.Range("K131").Value = 2607
$ws.Range("M131").Value = 2433

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2371.7334
$ws.Range("I132").Value = 1968.1538
$ws.Range("K132").Value = 5904.4614
$ws.Range("M132").Value = -3374.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4477.2856
$ws.Range("I7").Value = 2878.8333
$ws.Range("J7").Value = 5676.125
$ws.Range("K7").Value = 2878.8333
$ws.Range("L7").Value = 5676.125
$ws.Range("M7").Value = -2766.8333
$ws.Range("N7").Value = -5900.125

$ws.Range("H46").Value = 1211.75
$ws.Range("J46").Value = 1382
$ws.Range("L46").Value = 1382
$ws.Range("N46").Value = -1758

$ws.Range("H68").Value = 6185.7144
$ws.Range("I68").Value = 1615
$ws.Range("K68").Value = 1615
$ws.Range("M68").Value = -866

$ws.Range("H71").Value = 6185.7144
$ws.Range("I71").Value = 1615
$ws.Range("K71").Value = 8075
$ws.Range("M71").Value = -4331

$ws.Range("H122").Value = 6479.6665
$ws.Range("I122").Value = 4881.6816
$ws.Range("K122").Value = 14645.0448
$ws.Range("M122").Value = -12195.0448

$ws.Range("H126").Value = 4477.2856
$ws.Range("I126").Value = 2878.8333
$ws.Range("J126").Value = 5676.125
$ws.Range("K126").Value = 8636.499899999999
$ws.Range("L126").Value = 17028.375
$ws.Range("M126").Value = -6166.499899999999
$ws.Range("N126").Value = -21968.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3776.5334
$ws.Range("I122").Value = 3748.394
$ws.Range("J122").Value = 3853.9167
$ws.Range("K122").Value = 11245.182
$ws.Range("L122").Value = 11561.7501
$ws.Range("M122").Value = -8795.181999999999
$ws.Range("N122").Value = -16461.7501

$ws.Range("H132").Value = 1354.2
$ws.Range("I132").Value = 1338.2778
$ws.Range("K132").Value = 4014.8334
$ws.Range("M132").Value = -1484.8334

$ws.Range("H136").Value = 4601.8
$ws.Range("I136").Value = 2388.238
$ws.Range("K136").Value = 7164.714
$ws.Range("M136").Value = -4614.714

Write-Output "Applied 41 row updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"